# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing "Late", "Outstanding"
# and "Original" columns one place to the right, and updates the
# worksheet's active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column to the left (M) so the newly
# inserted column inherits the same width, matching Excel's default
# "Insert" behaviour of copying formatting from the column to its left.
$leftColumnWidth = $ws.Columns("M").ColumnWidth

# Insert a new column before column N - this shifts the old N/O/P
# columns (Late / Original / Outstanding) one column to the right.
$ws.Columns("N").Insert()

# Apply the inherited column width to the newly inserted column.
$ws.Columns("N").ColumnWidth = $leftColumnWidth

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("K14").Select()
